$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115-117 down to 116-118.
$ws.Rows.Item(115).Insert()

# The new row 115 repeats the same record pattern as the surrounding rows,
# just with different Fecha / Volumen / Precios / Precio $/Kg values.
$ws.Cells.Item(115, 1).Value = 9
$ws.Cells.Item(115, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(115, 3).Value = "Metropolitana"
$ws.Cells.Item(115, 4).Value = 44509
$ws.Cells.Item(115, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat
$ws.Cells.Item(115, 5).Value = 13
$ws.Cells.Item(115, 6).Value = "Fruta"
$ws.Cells.Item(115, 7).Value = 100101
$ws.Cells.Item(115, 8).Value = "Berries"
$ws.Cells.Item(115, 9).Value = 100101001
$ws.Cells.Item(115, 10).Value = "Arándano (blue)"
$ws.Cells.Item(115, 11).Value = "Sin especificar"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 830
$ws.Cells.Item(115, 14).Value = 8000
$ws.Cells.Item(115, 15).Value = 9000
$ws.Cells.Item(115, 16).Value = 8542
$ws.Cells.Item(115, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(115, 18).Value = "Provincia de Linares"
$ws.Cells.Item(115, 19).Value = 4271
$ws.Cells.Item(115, 20).Value = 2
